# "Revert to see more variables"
#
# The sheet currently has a blank spacer row (row 22) followed by the
# "avg_charges_tot" row (23), an "avgChronCond" row (24), and the
# "borocodenum" row (25). The edit removes the blank spacer row and the
# "avgChronCond" row entirely, shifting the remaining rows up so the sheet
# ends with "avg_charges_tot" at row 22 and "borocodenum" at row 23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank spacer row (row 22). Everything below shifts up by one:
# the old row 23 (avg_charges_tot) becomes row 22, old row 24 (avgChronCond)
# becomes row 23, old row 25 (borocodenum) becomes row 24.
$ws.Rows.Item(22).Delete()

# Delete the now-row-23 "avgChronCond" row. Old row 25 (borocodenum) shifts
# up to become row 23.
$ws.Rows.Item(23).Delete()

# Update the view to match: scroll near the bottom of the now-shorter list
# and select cell A20.
$ws.Range("A20").Select()
$excel.ActiveWindow.ScrollRow = 19
